$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New emotional-response values for columns B (Sad), C (Happy), D (Flirtatious), E (Anger)
# rows 2-21. Empty string clears a cell (keeps it blank like before).
$data = @(
    @("mournful",      "cake",      "bae",        "suck"),
    @("dismay",         "like",      "attractive", "bitch"),
    @("heartbreaking",  "puppies",   "attracted",  "ass"),
    @("pathetic",       "pizza",     "cute",       "fatty"),
    @("wound",          ":)",        ";)",         "fuck"),
    @("cheating",       "sunflower", "cutie",      "faggot"),
    @("miserable",      "cats",      "sexy",       "jerk"),
    @("stab",           "friends",   "babe",       "dork"),
    @("unhappy",        "success",   "crush",      "dumb"),
    @("sad",            "passed",    "sex",        "schmuck"),
    @("hang over",      "",          "kiss",       "hate"),
    @("death",          "",          "baby",       "cunt"),
    @("tragic",         "",          "love",       "bonehead"),
    @("terrible",       "",          "honey",      "liptard"),
    @("worry",          "",          "",           "stupid"),
    @("devastate",      "",          "",           "loser"),
    @("distress",       "",          "",           "fat"),
    @("disturb",        "",          "",           "lazy"),
    @("upset",          "",          "",           ""),
    @("depress",        "",          "",           "")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $vals = $data[$i]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
}
